$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52 <- values previously in row 54
$ws.Range("B52").Value = 6834770
$ws.Range('E52').Value = 'NK Solin'
$ws.Range('F52').Value = 'NK Croatia Zmijavci'
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range('K52').Value = 'D'
$ws.Range("L52").Value = 1.85
$ws.Range("M52").Value = 3.4
$ws.Range("N52").Value = 3.6
$ws.Range("O52").Value = 1.909
$ws.Range("P52").Value = 3.4
$ws.Range("Q52").Value = 3.5
$ws.Range("R52").Value = -0.5
$ws.Range("S52").Value = 1.95
$ws.Range("T52").Value = 1.85
$ws.Range("U52").Value = 2.5
$ws.Range("V52").Value = 1.9
$ws.Range("W52").Value = 1.9
$ws.Range("X52").Value = -1
$ws.Range("Y52").Value = 2.4
$ws.Range("Z52").Value = -1
$ws.Range("AA52").Value = -1
$ws.Range("AB52").Value = 0.8500000000000001
$ws.Range("AC52").Value = -1
$ws.Range("AD52").Value = 0.8999999999999999

# Row 54 <- values previously in row 52
$ws.Range("B54").Value = 6834769
$ws.Range('E54').Value = 'Vukovar 91'
$ws.Range('F54').Value = 'NK Dubrava Zagreb'
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 2
$ws.Range("I54").Value = 1
$ws.Range("J54").Value = 1
$ws.Range('K54').Value = 'D'
$ws.Range("L54").Value = 1.727
$ws.Range("M54").Value = 3.6
$ws.Range("N54").Value = 3.8
$ws.Range("O54").Value = 1.909
$ws.Range("P54").Value = 3.5
$ws.Range("Q54").Value = 3.2
$ws.Range("R54").Value = -0.5
$ws.Range("S54").Value = 2
$ws.Range("T54").Value = 1.8
$ws.Range("U54").Value = 2.5
$ws.Range("V54").Value = 1.975
$ws.Range("W54").Value = 1.825
$ws.Range("X54").Value = -1
$ws.Range("Y54").Value = 2.5
$ws.Range("Z54").Value = -1
$ws.Range("AA54").Value = -1
$ws.Range("AB54").Value = 0.8
$ws.Range("AC54").Value = 0.9750000000000001
$ws.Range("AD54").Value = -1

# Row 147 <- values previously in row 148
$ws.Range("B147").Value = 7977245
$ws.Range('E147').Value = 'NK Dugopolje'
$ws.Range('F147').Value = 'NK Dubrava Zagreb'
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 0
$ws.Range("I147").Value = 0
$ws.Range("J147").Value = 0
$ws.Range('K147').Value = 'D'
$ws.Range("L147").Value = 2
$ws.Range("M147").Value = 3.1
$ws.Range("N147").Value = 3.4
$ws.Range("O147").Value = 2.7
$ws.Range("P147").Value = 3
$ws.Range("Q147").Value = 2.55
$ws.Range("R147").Value = 0
$ws.Range("S147").Value = 1.95
$ws.Range("T147").Value = 1.85
$ws.Range("U147").Value = 2.25
$ws.Range("V147").Value = 1.925
$ws.Range("W147").Value = 1.875
$ws.Range("X147").Value = -1
$ws.Range("Y147").Value = 2
$ws.Range("Z147").Value = -1
$ws.Range("AA147").Value = 0
$ws.Range("AB147").Value = 0
$ws.Range("AC147").Value = -1
$ws.Range("AD147").Value = 0.875

# Row 148 <- values previously in row 147
$ws.Range("B148").Value = 7977243
$ws.Range('E148').Value = 'Bijelo Brdo'
$ws.Range('F148').Value = 'NK Sesvete'
$ws.Range("G148").Value = 2
$ws.Range("H148").Value = 1
$ws.Range("I148").Value = 0
$ws.Range("J148").Value = 1
$ws.Range('K148').Value = 'H'
$ws.Range("L148").Value = 2.4
$ws.Range("M148").Value = 3.4
$ws.Range("N148").Value = 2.5
$ws.Range("O148").Value = 1.909
$ws.Range("P148").Value = 3.2
$ws.Range("Q148").Value = 4
$ws.Range("R148").Value = -0.5
$ws.Range("S148").Value = 1.9
$ws.Range("T148").Value = 1.9
$ws.Range("U148").Value = 2.25
$ws.Range("V148").Value = 1.975
$ws.Range("W148").Value = 1.725
$ws.Range("X148").Value = 0.909
$ws.Range("Y148").Value = -1
$ws.Range("Z148").Value = -1
$ws.Range("AA148").Value = 0.8999999999999999
$ws.Range("AB148").Value = -1
$ws.Range("AC148").Value = 0.9750000000000001
$ws.Range("AD148").Value = -1

# Row 201 <- values previously in row 202
$ws.Range("B201").Value = 7977299
$ws.Range('E201').Value = 'Bijelo Brdo'
$ws.Range('F201').Value = 'NK Dugopolje'
$ws.Range("G201").Value = 2
$ws.Range("H201").Value = 1
$ws.Range('K201').Value = 'H'
$ws.Range("L201").Value = 1.363
$ws.Range("M201").Value = 4.333
$ws.Range("N201").Value = 7
$ws.Range("O201").Value = 1.4
$ws.Range("P201").Value = 4.333
$ws.Range("Q201").Value = 6
$ws.Range("R201").Value = -1.25
$ws.Range("S201").Value = 1.95
$ws.Range("T201").Value = 1.85
$ws.Range("U201").Value = 2.5
$ws.Range("V201").Value = 1.875
$ws.Range("W201").Value = 1.925
$ws.Range("X201").Value = 0.3999999999999999
$ws.Range("Y201").Value = -1
$ws.Range("Z201").Value = -1
$ws.Range("AA201").Value = -0.5
$ws.Range("AB201").Value = 0.425
$ws.Range("AC201").Value = 0.875
$ws.Range("AD201").Value = -1

# Row 202 <- values previously in row 203
$ws.Range("B202").Value = 7977300
$ws.Range('E202').Value = 'NK Dubrava Zagreb'
$ws.Range('F202').Value = 'NK Solin'
$ws.Range("G202").Value = 3
$ws.Range("H202").Value = 0
$ws.Range('K202').Value = 'H'
$ws.Range("L202").Value = 1.85
$ws.Range("M202").Value = 3.25
$ws.Range("N202").Value = 3.75
$ws.Range("O202").Value = 1.5
$ws.Range("P202").Value = 4.1
$ws.Range("Q202").Value = 5
$ws.Range("R202").Value = -1
$ws.Range("S202").Value = 1.85
$ws.Range("T202").Value = 1.95
$ws.Range("U202").Value = 3
$ws.Range("V202").Value = 1.9
$ws.Range("W202").Value = 1.9
$ws.Range("X202").Value = 0.5
$ws.Range("Y202").Value = -1
$ws.Range("Z202").Value = -1
$ws.Range("AA202").Value = 0.8500000000000001
$ws.Range("AB202").Value = -1
$ws.Range("AC202").Value = 0
$ws.Range("AD202").Value = 0

# Row 203 <- values previously in row 201
$ws.Range("B203").Value = 7977297
$ws.Range('E203').Value = 'NK Croatia Zmijavci'
$ws.Range('F203').Value = 'NK Sesvete'
$ws.Range("G203").Value = 3
$ws.Range("H203").Value = 1
$ws.Range('K203').Value = 'H'
$ws.Range("L203").Value = 1.5
$ws.Range("M203").Value = 3.4
$ws.Range("N203").Value = 6.5
$ws.Range("O203").Value = 1.4
$ws.Range("P203").Value = 3.8
$ws.Range("Q203").Value = 6.5
$ws.Range("R203").Value = -1.25
$ws.Range("S203").Value = 1.925
$ws.Range("T203").Value = 1.875
$ws.Range("U203").Value = 3
$ws.Range("V203").Value = 2
$ws.Range("W203").Value = 1.8
$ws.Range("X203").Value = 0.3999999999999999
$ws.Range("Y203").Value = -1
$ws.Range("Z203").Value = -1
$ws.Range("AA203").Value = 0.925
$ws.Range("AB203").Value = -1
$ws.Range("AC203").Value = 1
$ws.Range("AD203").Value = -1
